# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 42 (pushing the existing rows 42-103
# down to 43-104) and populate it with the new week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(42).Insert()

$ws.Range("A42").Value = 8
$ws.Range("B42").Value = "Terminal La Palmera de La Serena"
$ws.Range("C42").Value = "Coquimbo"
$ws.Range("D42").Value = 44540
$ws.Range("E42").Value = 4
$ws.Range("F42").Value = 100112040
$ws.Range("G42").Value = "Cilantro"
$ws.Range("H42").Value = "Sin especificar"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 3000
$ws.Range("K42").Value = 1500
$ws.Range("L42").Value = 2000
$ws.Range("M42").Value = 1750
$ws.Range("N42").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O42").Value = "Provincia del Elquí"
$ws.Range("P42").Value = 1167
$ws.Range("Q42").Value = 1.5
$ws.Range("R42").Value = "Hortaliza"
